$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57; existing rows 57-65 shift down to 58-66.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with its data.
$ws.Cells.Item(57, 1).Value = 10
$ws.Cells.Item(57, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(57, 3).Value = "La Araucanía"
$ws.Cells.Item(57, 4).Value = 44826
$ws.Cells.Item(57, 5).Value = 9
$ws.Cells.Item(57, 6).Value = 300000000
$ws.Cells.Item(57, 7).Value = "Espárragos"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 100
$ws.Cells.Item(57, 11).Value = 3000
$ws.Cells.Item(57, 12).Value = 3000
$ws.Cells.Item(57, 13).Value = 3000
$ws.Cells.Item(57, 14).Value = "$/kilo"
$ws.Cells.Item(57, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(57, 16).Value = 3000
$ws.Cells.Item(57, 17).Value = 1
$ws.Cells.Item(57, 18).Value = "Hortaliza"
